# Updates cryptos list Price (D) and Volume(1h) (E) columns per upstream refresh.
# Price/percentage cells are plain text in this sheet (no numeric semantics),
# so for Price values that would otherwise auto-parse as a number we pin the
# cell's NumberFormat to Text ("@") first, exactly as a user re-typing over an
# existing text cell would end up doing, to keep them stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.983.98"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "2.609.45"
$ws.Range("E3").Value = "  -2.28%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.51"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.18"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("D9").Value = "2.608.78"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.21"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.48"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Value = "3.087.90"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "66.935.98"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").Value = "2.641.68"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.97"
$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  +5.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.11"
$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.63"
$ws.Range("E23").Value = "  -3.28%  "

$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("E25").Value = "  -5.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.15"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("E27").Value = "  -2.09%  "

$ws.Range("D28").Value = "2.749.97"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "0.0₃0994"
$ws.Range("E30").Value = "  -2.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "540.76"
$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  -1.49%  "

$ws.Range("E33").Value = "  -2.99%  "

$ws.Range("E34").Value = "  -1.62%  "

$ws.Range("E35").Value = "  +4.91%  "

$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("E37").Value = "  -4.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.80"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.91"
$ws.Range("E39").Value = "  -2.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.364"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -4.67%  "

$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.28"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.574"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("E51").Value = "  -0.64%  "
